$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old column A (index numbers) and column B (employee names)
# so the old strings ("Row 1", "Jack", "Martin", ...) are dropped and the
# new ones are interned fresh, in entry order.
$ws.Range("A1:B8").ClearContents()

# Column A now holds "Name" header + the employee names (was a plain numeric index 0..6).
$names = @("Name", "Dirk", "Anna", "Marie", "John", "Peter", "Stephanie", "Laura")
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $names[$i]
}

# Column B now holds the question header (unchanged text) + the numeric ratings
# (was the employee name text).
$ws.Cells.Item(1, 2).Value = "How motivated are you to come to work every day?"
$values = @(7, 6, 8, 3, 9, 6, 4)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

# Update the selection to match the saved view state.
$ws.Range("A2:A8").Select()
$excel.ActiveCell = $ws.Range("A8")

# The saved file now also carries an explicit page setup (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
